$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that Excel would otherwise try to "smart parse" (like a
# date-shaped string such as "2017-04-12") as plain text, exactly as typed, and
# without leaving behind any extra cell formatting/style.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.ClearFormats()
}

# New "Comment" column header (H1)
$ws.Range("H1").Value = "Comment"

# Row 2 (Coach 10805 / Axle 4 / Wheel 8)
$ws.Range("F2").Value = 5
Set-TextValue $ws.Range("G2") "2017-04-12"
$ws.Range("H2").Value = "The Toe Creep has violated the wheel alarm settings thresholds"

# Row 3 (Coach 10805 / Axle 4 / Wheel 7)
$ws.Range("F3").Value = 4
Set-TextValue $ws.Range("G3") "2017-04-11"
$ws.Range("H3").Value = "The Toe Creep has violated the wheel alarm settings thresholds"

# Row 4 (Coach 10805 / Axle 3 / Wheel 6)
$ws.Range("F4").Value = 5
Set-TextValue $ws.Range("G4") "2017-04-12"
$ws.Range("H4").Value = "The Toe Creep has violated the wheel alarm settings thresholds"

# Row 5 (Coach 10805 / Axle 3 / Wheel 5)
$ws.Range("F5").Value = 5
Set-TextValue $ws.Range("G5") "2017-04-12"
$ws.Range("H5").Value = "The Toe Creep has violated the wheel alarm settings thresholds"

# Row 6 (Coach 10805 / Axle 2 / Wheel 4) - F6/G6 unchanged, only new Comment added
$ws.Range("H6").Value = "The Flange Height has violated the wheel alarm settings thresholds"
